$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'45.354.74"
$ws.Range("E2").Value = "  +2.77%  "
$ws.Range("D3").Value = "'2.425.36"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'318.20"
$ws.Range("E5").Value = "  +3.30%  "
$ws.Range("D6").Value = "'102.57"
$ws.Range("E6").Value = "  +4.62%  "
$ws.Range("D7").Value = "'0.516"
$ws.Range("E7").Value = "  +0.85%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.529"
$ws.Range("E9").Value = "  +6.31%  "
$ws.Range("D10").Value = "'35.58"
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("D11").Value = "'0.0803"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("E12").Value = "  -2.16%  "
$ws.Range("D13").Value = "'18.11"
$ws.Range("E13").Value = "  -2.99%  "
$ws.Range("D14").Value = "'7.06"
$ws.Range("E14").Value = "  +1.68%  "
$ws.Range("D15").Value = "'2.807.66"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").Value = "'2.431.46"
$ws.Range("E16").Value = "  -2.40%  "
$ws.Range("D17").Value = "'0.843"
$ws.Range("E17").Value = "  +1.36%  "
$ws.Range("D18").Value = "'45.317.42"
$ws.Range("E18").Value = "  +2.80%  "
$ws.Range("D19").Value = "'12.24"
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").Value = "'6.34"
$ws.Range("E20").Value = "  -1.84%  "
$ws.Range("D21").Value = "'0.0₃0921"
$ws.Range("E21").Value = "  +1.71%  "
$ws.Range("D22").Value = "'68.83"
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("D23").Value = "'244.14"
$ws.Range("E23").Value = "  +1.88%  "
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("D25").Value = "'2.49"
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Value = "'25.54"
$ws.Range("E27").Value = "  +1.57%  "
$ws.Range("E28").Value = "  -1.85%  "
$ws.Range("D29").Value = "'9.57"
$ws.Range("E29").Value = "  +1.01%  "
$ws.Range("D30").Value = "'49.16"
$ws.Range("E30").Value = "  +2.47%  "
$ws.Range("D31").Value = "'32.75"
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("D32").Value = "'20.32"
$ws.Range("E32").Value = "  +9.25%  "
$ws.Range("E33").Value = "  +4.32%  "
$ws.Range("D34").Value = "'5.21"
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("E35").Value = "  +0.43%  "
$ws.Range("D36").Value = "'0.0767"
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("D37").Value = "'1.87"
$ws.Range("E37").Value = "  -2.84%  "
$ws.Range("D38").Value = "'4.44"
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("E39").Value = "  -2.15%  "
$ws.Range("D40").Value = "'125.63"
$ws.Range("E40").Value = "  -4.36%  "
$ws.Range("E41").Value = "  -2.53%  "
$ws.Range("D42").Value = "'0.109"
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("D43").Value = "'20.48"
$ws.Range("E43").Value = "  -4.98%  "
$ws.Range("E44").Value = "  +1.39%  "
$ws.Range("D45").Value = "'1.924.48"
$ws.Range("E45").Value = "  -1.49%  "
$ws.Range("E46").Value = "  -2.87%  "
$ws.Range("D47").Value = "'2.92"
$ws.Range("E47").Value = "  +1.72%  "
$ws.Range("D48").Value = "'1.81"
$ws.Range("E48").Value = "  +12.67%  "
$ws.Range("E49").Value = "  -2.19%  "
$ws.Range("D50").Value = "'76.58"
$ws.Range("E50").Value = "  +5.00%  "
$ws.Range("D51").Value = "'53.92"
$ws.Range("E51").Value = "  +1.51%  "
